$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Total "VALOR MORA" figure near the top of the account-statement header.
$ws.Range("E11").Value = 64940

# Replace the single worker row (row 16) with the new worker's data:
#  - Tipo Doc Trabajador : NIT -> CC
#  - N° Doc Trabajador   : 9007300963 -> 1073810296
#  - Nombre Trabajador   : (blank) -> MANUEL JOSE BANDA MARTINEZ
#  - Periodo Mora        : 1606 -> 2508
#  - Valor Mora          : 485000 -> 64940
#  - Salario Basico      : 0 -> 1623500
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1073810296"
$ws.Range("D16").Value = "MANUEL JOSE BANDA MARTINEZ"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 64940
$ws.Range("G16").Value = 1623500

# Column D ("Nombre Trabajador") is a best-fit column; widen it so the new,
# longer worker name is not clipped (mirrors Excel auto-resizing the column
# when the underlying cell content grows).
$ws.Columns.Item(4).AutoFit() | Out-Null
